$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing key/value pairs (A1:B7) before rearranging them.
$headers = @()
$values = @()
for ($r = 1; $r -le 7; $r++) {
    $headers += $ws.Cells.Item($r, 1).Value2
    $values += $ws.Cells.Item($r, 2).Value2
}

# Clear the old layout.
$ws.Range("A1:B7").Clear()

# Rebuild as a transposed table: row 1 = headers, row 2 = values,
# spread across columns A:G.
for ($i = 0; $i -lt 7; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
